$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name and title text to reflect new "through" date
$ws.Name = "Through 2022-07-02"
$ws.Range("A8").Value = "July (through 07-02)"

# Update July row (row 8) values
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 9
$ws.Range("I8").Value = 15

# Update Total row (row 9) values
$ws.Range("C9").Value = 252
$ws.Range("D9").Value = 394
$ws.Range("E9").Value = 360
$ws.Range("F9").Value = 254
$ws.Range("G9").Value = 475
$ws.Range("H9").Value = 769
$ws.Range("I9").Value = 821
